$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.175.21"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "2.052.55"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.95"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("E6").Value = "  +1.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.36"
$ws.Range("E7").Value = "  +6.67%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +0.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0785"
$ws.Range("E10").Value = "  -1.28%  "
$ws.Range("E11").Value = "  +1.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.20"
$ws.Range("E12").Value = "  +6.69%  "
$ws.Range("D13").Value = "2.352.35"
$ws.Range("E13").Value = "  -0.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.817"
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("E15").Value = "  +6.42%  "
$ws.Range("D16").Value = "2.050.52"
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("D17").Value = "37.140.54"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.75"
$ws.Range("E18").Value = "  +24.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.90"
$ws.Range("E19").Value = "  +3.40%  "
$ws.Range("E20").Value = "  -1.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.39"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "237.60"
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("E24").Value = "  -0.94%  "
$ws.Range("E25").Value = "  +11.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.59"
$ws.Range("E26").Value = "  -0.88%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.36"
$ws.Range("E27").Value = "  +3.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.95"
$ws.Range("E28").Value = "  -1.26%  "
$ws.Range("E29").Value = "  +1.32%  "
$ws.Range("E30").Value = "  +8.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.76"
$ws.Range("E31").Value = "  +4.61%  "
$ws.Range("E32").Value = "  -1.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.51"
$ws.Range("E33").Value = "  +3.97%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0896"
$ws.Range("E34").Value = "  +3.75%  "
$ws.Range("E35").Value = "  -0.14%  "
$ws.Range("E36").Value = "  -1.24%  "
$ws.Range("E37").Value = "  -1.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.111"
$ws.Range("E38").Value = "  +7.71%  "
$ws.Range("E39").Value = "  +0.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.30"
$ws.Range("E40").Value = "  +31.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.15"
$ws.Range("E41").Value = "  +12.84%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.71"
$ws.Range("E42").Value = "  -2.03%  "
$ws.Range("E43").Value = "  -0.94%  "
$ws.Range("E44").Value = "  -1.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "96.29"
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("E46").Value = "  +1.90%  "
$ws.Range("D47").Value = "1.286.19"
$ws.Range("E47").Value = "  -0.94%  "
$ws.Range("E48").Value = "  -1.36%  "
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("D50").Value = "2.244.21"
$ws.Range("E50").Value = "  -0.37%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "42.85"
$ws.Range("E51").Value = "  -2.28%  "
